$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = "我的收入 - 2.psd"
$ws.Range("A84").Value = "收入明细.psd"

$ws.Range("B82").Value = "income.html"
$ws.Range("B84").Value = "income-detail.html"

$ws.Range("A85").Value = "收入明细-详情页面.psd"
$ws.Range("A86").Value = "收入明细-详情页面个人明细和扩展团队.psd"

$ws.Range("B85").Value = "income-detail-info.html"
$ws.Range("B86").Value = "income-detail-info2.html"

[void]$ws.Range("A76").Select()
